$wb = $excel.ActiveWorkbook

# --- Sheet "Era C" (sheet1) ---
$ws = $wb.Worksheets.Item("Era C")

$ws.Range("B2").Value = 57327.05
$ws.Range("C2").Value = 69089.91
$ws.Range("D2").Value = 126416.96

$ws.Range("B3").Value = -11891.18
$ws.Range("C3").Value = -20277.56
$ws.Range("D3").Value = -32168.74

$ws.Range("B4").Value = -28564.63
$ws.Range("C4").Value = -29994.23
$ws.Range("D4").Value = -58558.86

$ws.Range("B5").Value = -196582.08
$ws.Range("C5").Value = -188143.44
$ws.Range("D5").Value = -384725.52

$ws.Range("B6").Value = -18135.67
$ws.Range("C6").Value = -16482.58
$ws.Range("D6").Value = -34618.25

$ws.Range("B7").Value = 81716.10000000001
$ws.Range("C7").Value = 76255.21000000001
$ws.Range("D7").Value = 157971.31

$ws.Range("B8").Value = -116130.41
$ws.Range("C8").Value = -109552.69
$ws.Range("D8").Value = -225683.1

# --- Sheet "Era B" (sheet2) ---
$ws = $wb.Worksheets.Item("Era B")

$ws.Range("B2").Value = 46505.97
$ws.Range("C2").Value = 46505.97

$ws.Range("B3").Value = -16132.08
$ws.Range("C3").Value = -16132.08

$ws.Range("B4").Value = -27644.43
$ws.Range("C4").Value = -27644.43

$ws.Range("B5").Value = -185434.74
$ws.Range("C5").Value = -185434.74

$ws.Range("B6").Value = -15629.56
$ws.Range("C6").Value = -15629.56

$ws.Range("B7").Value = 76187.89999999999
$ws.Range("C7").Value = 76187.89999999999

$ws.Range("B8").Value = -122146.94
$ws.Range("C8").Value = -122146.94

# --- Sheet "Era A" (sheet3) ---
$ws = $wb.Worksheets.Item("Era A")

$ws.Range("B2").Value = 68042.99000000001
$ws.Range("C2").Value = 97647.47
$ws.Range("D2").Value = 165690.46

$ws.Range("B3").Value = -18620.19
$ws.Range("C3").Value = -13311.51
$ws.Range("D3").Value = -31931.7

$ws.Range("B4").Value = -29930.08
$ws.Range("C4").Value = -24965.85
$ws.Range("D4").Value = -54895.93

$ws.Range("B5").Value = -205682.74
$ws.Range("C5").Value = -189920.85
$ws.Range("D5").Value = -395603.59

$ws.Range("B6").Value = -19907.67
$ws.Range("C6").Value = -12846.96
$ws.Range("D6").Value = -32754.63

$ws.Range("B7").Value = 68580.99000000001
$ws.Range("C7").Value = 110584.34
$ws.Range("D7").Value = 179165.33

$ws.Range("B8").Value = -137516.7
$ws.Range("C8").Value = -32813.36
$ws.Range("D8").Value = -170330.06
